# Apply NATMI LR-pairs update: add M2 sending-cluster rows and refresh
# statistics for the existing M1 sending-cluster rows (rows 2-6),
# then append the new M2 sending-cluster rows (rows 7-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowData = @{}
$rowData[2] = @("M1", "Cd177", "Pecam1", "ECs", 3, 1, 2.166130333333333, 6.498391, 0.237044430060999, 0.237044430060999, 3, 1, 328.5062863333333, 985.518859, 0.9485223498994318, 0.9514072518130269, 711.5874315173187, 6404.286883655869, 0.2248419398320302, 0.2255257897619203)
$rowData[3] = @("M1", "Cd177", "Pecam1", "FAPs", 3, 1, 2.166130333333333, 6.498391, 0.237044430060999, 0.237044430060999, 3, 1, 4.781633, 14.344899, 0.01380638958279946, 0.0138483812973136, 10.35764028416767, 93.218762557509, 0.003272727749854812, 0.0032826816518891)
$rowData[4] = @("M1", "Cd177", "Pecam1", "M1", 3, 1, 2.166130333333333, 6.498391, 0.237044430060999, 0.237044430060999, 3, 1, 5.094370666666666, 15.283112, 0.0147093819419403, 0.01475412008028422, 11.03507083031022, 99.31563747279199, 0.003486777058976788, 0.003497381985482513)
$rowData[5] = @("M1", "Cd177", "Pecam1", "M2", 3, 1, 2.166130333333333, 6.498391, 0.237044430060999, 0.237044430060999, 3, 1, 4.801979333333333, 14.405938, 0.01386513717061757, 0.01390730756413546, 10.40171309397311, 93.61541784575799, 0.003286653538326614, 0.003296649795223511)
$rowData[6] = @("M1", "Cd177", "Pecam1", "sCs", 3, 1, 2.166130333333333, 6.498391, 0.237044430060999, 0.237044430060999, 2, 1, 3.150518, 6.301036, 0.00909674140521077, 0.006082939245239696, 6.824432605512667, 40.946595633076, 0.002156331881810478, 0.001441926866483527)
$rowData[7] = @("M2", "Cd177", "Pecam1", "ECs", 2, 0.6666666666666666, 6.971947, 20.915841, 0.762955569939001, 0.762955569939001, 3, 1, 328.5062863333333, 985.518859, 0.9485223498994318, 0.9514072518130269, 2290.328417482824, 20612.95575734542, 0.7236804100674016, 0.7258814620511066)
$rowData[8] = @("M2", "Cd177", "Pecam1", "FAPs", 2, 0.6666666666666666, 6.971947, 20.915841, 0.762955569939001, 0.762955569939001, 3, 1, 4.781633, 14.344899, 0.01380638958279946, 0.0138483812973136, 33.337291849451, 300.035626645059, 0.01053366183294465, 0.0105656996454245)
$rowData[9] = @("M2", "Cd177", "Pecam1", "M1", 2, 0.6666666666666666, 6.971947, 20.915841, 0.762955569939001, 0.762955569939001, 3, 1, 5.094370666666666, 15.283112, 0.0147093819419403, 0.01475412008028422, 35.51768228635466, 319.659140577192, 0.01122260488296351, 0.01125673809480171)
$rowData[10] = @("M2", "Cd177", "Pecam1", "M2", 2, 0.6666666666666666, 6.971947, 20.915841, 0.762955569939001, 0.762955569939001, 3, 1, 4.801979333333333, 14.405938, 0.01386513717061757, 0.01390730756413546, 33.47914540709533, 301.312308663858, 0.01057848363229096, 0.01061065776891195)
$rowData[11] = @("M2", "Cd177", "Pecam1", "sCs", 2, 0.6666666666666666, 6.971947, 20.915841, 0.762955569939001, 0.762955569939001, 2, 1, 3.150518, 6.301036, 0.00909674140521077, 0.006082939245239696, 21.965244518546, 131.791467111276, 0.006940409523400292, 0.00464101237875617)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in 2..11) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

Write-Host ("New dimension: " + $ws.UsedRange.Address())
